# Data base update: reorder the "Periodo Mora" (year) rows 16-22 so the
# years run newest-to-oldest (2006 down to 1912), carry the "Valor Mora"
# figures along with their year, and zero out "Salario Basico" for all
# these rows (previous account-statement periods were cleared / new ones
# added, per the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$years  = @("2006", "2005", "2004", "2003", "2002", "2001", "1912")
$valorMora = @(26500, 33125, 33125, 33125, 33125, 33125, 33125)
$salarioBasico = @(0, 0, 0, 0, 0, 0, 0)

$startRow = 16
for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $years[$i]
    $ws.Cells.Item($row, 6).Value = $valorMora[$i]
    $ws.Cells.Item($row, 7).Value = $salarioBasico[$i]
}
